# Updates the cryptos price table (Price / Volume(1h) columns, plus a row
# reorder for Hedera/VeChain) to reflect the refreshed GitHub Actions data.
#
# Several "Price" values are plain decimal-looking strings (e.g. "1.006",
# "0.00001044") that Excel would otherwise auto-convert to numbers (losing
# trailing zeros / switching to scientific notation) if assigned directly
# via .Value. For those cells we briefly force text format ("@"), assign
# the literal string, then ClearFormats() so the cell's style reverts to
# the workbook default (matching the original, unstyled data cells) while
# the text content itself is preserved exactly as authored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = '24.551.42'
$ws.Cells.Item(2,5).Value = '  -1.24%  '
$ws.Cells.Item(3,4).Value = '1.654.92'
$ws.Cells.Item(3,5).Value = '  -2.94%  '
$ws.Cells.Item(4,4).NumberFormat = '@'
$ws.Cells.Item(4,4).Value = '1.006'
$ws.Cells.Item(4,4).ClearFormats()
$ws.Cells.Item(4,5).Value = '  -0.02%  '
$ws.Cells.Item(5,4).NumberFormat = '@'
$ws.Cells.Item(5,4).Value = '318.92'
$ws.Cells.Item(5,4).ClearFormats()
$ws.Cells.Item(5,5).Value = '  +2.11%  '
$ws.Cells.Item(6,4).NumberFormat = '@'
$ws.Cells.Item(6,4).Value = '1.001'
$ws.Cells.Item(6,4).ClearFormats()
$ws.Cells.Item(6,5).Value = '  +0.11%  '
$ws.Cells.Item(7,4).NumberFormat = '@'
$ws.Cells.Item(7,4).Value = '0.3634'
$ws.Cells.Item(7,4).ClearFormats()
$ws.Cells.Item(7,5).Value = '  -2.88%  '
$ws.Cells.Item(8,4).NumberFormat = '@'
$ws.Cells.Item(8,4).Value = '46.83'
$ws.Cells.Item(8,4).ClearFormats()
$ws.Cells.Item(8,5).Value = '  -5.13%  '
$ws.Cells.Item(9,4).NumberFormat = '@'
$ws.Cells.Item(9,4).Value = '0.3257'
$ws.Cells.Item(9,4).ClearFormats()
$ws.Cells.Item(9,5).Value = '  -5.11%  '
$ws.Cells.Item(10,4).NumberFormat = '@'
$ws.Cells.Item(10,4).Value = '1.132'
$ws.Cells.Item(10,4).ClearFormats()
$ws.Cells.Item(10,5).Value = '  -7.28%  '
$ws.Cells.Item(11,4).NumberFormat = '@'
$ws.Cells.Item(11,4).Value = '0.07042'
$ws.Cells.Item(11,4).ClearFormats()
$ws.Cells.Item(11,5).Value = '  -6.23%  '
$ws.Cells.Item(12,4).NumberFormat = '@'
$ws.Cells.Item(12,4).Value = '1.000'
$ws.Cells.Item(12,4).ClearFormats()
$ws.Cells.Item(12,5).Value = '  -0.14%  '
$ws.Cells.Item(13,4).NumberFormat = '@'
$ws.Cells.Item(13,4).Value = '6.028'
$ws.Cells.Item(13,4).ClearFormats()
$ws.Cells.Item(13,5).Value = '  -4.98%  '
$ws.Cells.Item(14,4).NumberFormat = '@'
$ws.Cells.Item(14,4).Value = '19.49'
$ws.Cells.Item(14,4).ClearFormats()
$ws.Cells.Item(14,5).Value = '  -7.74%  '
$ws.Cells.Item(15,4).Value = '1.660.97'
$ws.Cells.Item(15,5).Value = '  -2.79%  '
$ws.Cells.Item(16,4).NumberFormat = '@'
$ws.Cells.Item(16,4).Value = '6.613'
$ws.Cells.Item(16,4).ClearFormats()
$ws.Cells.Item(16,5).Value = '  -6.04%  '
$ws.Cells.Item(17,4).NumberFormat = '@'
$ws.Cells.Item(17,4).Value = '0.00001044'
$ws.Cells.Item(17,4).ClearFormats()
$ws.Cells.Item(17,5).Value = '  -7.72%  '
$ws.Cells.Item(18,4).NumberFormat = '@'
$ws.Cells.Item(18,4).Value = '0.06626'
$ws.Cells.Item(18,4).ClearFormats()
$ws.Cells.Item(18,5).Value = '  -1.42%  '
$ws.Cells.Item(19,4).NumberFormat = '@'
$ws.Cells.Item(19,4).Value = '0.9995'
$ws.Cells.Item(19,4).ClearFormats()
$ws.Cells.Item(19,5).Value = '  +0.03%  '
$ws.Cells.Item(20,4).NumberFormat = '@'
$ws.Cells.Item(20,4).Value = '78.85'
$ws.Cells.Item(20,4).ClearFormats()
$ws.Cells.Item(20,5).Value = '  -5.92%  '
$ws.Cells.Item(21,4).NumberFormat = '@'
$ws.Cells.Item(21,4).Value = '5.923'
$ws.Cells.Item(21,4).ClearFormats()
$ws.Cells.Item(21,5).Value = '  -7.00%  '
$ws.Cells.Item(22,4).NumberFormat = '@'
$ws.Cells.Item(22,4).Value = '15.72'
$ws.Cells.Item(22,4).ClearFormats()
$ws.Cells.Item(22,5).Value = '  -9.41%  '
$ws.Cells.Item(23,4).NumberFormat = '@'
$ws.Cells.Item(23,4).Value = '12.54'
$ws.Cells.Item(23,4).ClearFormats()
$ws.Cells.Item(23,5).Value = '  -4.88%  '
$ws.Cells.Item(24,4).Value = '24.581.85'
$ws.Cells.Item(24,5).Value = '  -1.20%  '
$ws.Cells.Item(25,4).NumberFormat = '@'
$ws.Cells.Item(25,4).Value = '2.453'
$ws.Cells.Item(25,4).ClearFormats()
$ws.Cells.Item(25,5).Value = '  +0.24%  '
$ws.Cells.Item(26,4).NumberFormat = '@'
$ws.Cells.Item(26,4).Value = '2.380'
$ws.Cells.Item(26,4).ClearFormats()
$ws.Cells.Item(26,5).Value = '  -14.59%  '
$ws.Cells.Item(27,4).NumberFormat = '@'
$ws.Cells.Item(27,4).Value = '147.85'
$ws.Cells.Item(27,4).ClearFormats()
$ws.Cells.Item(27,5).Value = '  -1.22%  '
$ws.Cells.Item(28,4).NumberFormat = '@'
$ws.Cells.Item(28,4).Value = '18.62'
$ws.Cells.Item(28,4).ClearFormats()
$ws.Cells.Item(28,5).Value = '  -8.71%  '
$ws.Cells.Item(29,4).Value = '1.842.66'
$ws.Cells.Item(29,5).Value = '  -2.83%  '
$ws.Cells.Item(30,4).NumberFormat = '@'
$ws.Cells.Item(30,4).Value = '1.205'
$ws.Cells.Item(30,4).ClearFormats()
$ws.Cells.Item(30,5).Value = '  -4.23%  '
$ws.Cells.Item(31,4).NumberFormat = '@'
$ws.Cells.Item(31,4).Value = '125.80'
$ws.Cells.Item(31,4).ClearFormats()
$ws.Cells.Item(31,5).Value = '  -5.22%  '
$ws.Cells.Item(32,4).NumberFormat = '@'
$ws.Cells.Item(32,4).Value = '4.063'
$ws.Cells.Item(32,4).ClearFormats()
$ws.Cells.Item(32,5).Value = '  -3.88%  '
$ws.Cells.Item(33,4).NumberFormat = '@'
$ws.Cells.Item(33,4).Value = '5.830'
$ws.Cells.Item(33,4).ClearFormats()
$ws.Cells.Item(33,5).Value = '  -14.16%  '
$ws.Cells.Item(34,4).NumberFormat = '@'
$ws.Cells.Item(34,4).Value = '0.08464'
$ws.Cells.Item(34,4).ClearFormats()
$ws.Cells.Item(34,5).Value = '  -3.49%  '
$ws.Cells.Item(35,4).NumberFormat = '@'
$ws.Cells.Item(35,4).Value = '1.674'
$ws.Cells.Item(35,4).ClearFormats()
$ws.Cells.Item(35,5).Value = '  -5.61%  '
$ws.Cells.Item(36,4).NumberFormat = '@'
$ws.Cells.Item(36,4).Value = '12.27'
$ws.Cells.Item(36,4).ClearFormats()
$ws.Cells.Item(36,5).Value = '  -10.82%  '
$ws.Cells.Item(37,4).NumberFormat = '@'
$ws.Cells.Item(37,4).Value = '1.276'
$ws.Cells.Item(37,4).ClearFormats()
$ws.Cells.Item(37,5).Value = '  +0.22%  '
$ws.Cells.Item(38,4).NumberFormat = '@'
$ws.Cells.Item(38,4).Value = '5.207'
$ws.Cells.Item(38,4).ClearFormats()
$ws.Cells.Item(38,5).Value = '  -7.29%  '
$ws.Cells.Item(39,2).Value = 'Hedera'
$ws.Cells.Item(39,3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(39,4).NumberFormat = '@'
$ws.Cells.Item(39,4).Value = '0.06021'
$ws.Cells.Item(39,4).ClearFormats()
$ws.Cells.Item(39,5).Value = '  -9.30%  '
$ws.Cells.Item(40,2).Value = 'VeChain'
$ws.Cells.Item(40,3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(40,4).NumberFormat = '@'
$ws.Cells.Item(40,4).Value = '0.02229'
$ws.Cells.Item(40,4).ClearFormats()
$ws.Cells.Item(40,5).Value = '  -7.51%  '
$ws.Cells.Item(41,4).NumberFormat = '@'
$ws.Cells.Item(41,4).Value = '0.2071'
$ws.Cells.Item(41,4).ClearFormats()
$ws.Cells.Item(41,5).Value = '  -8.33%  '
$ws.Cells.Item(42,4).NumberFormat = '@'
$ws.Cells.Item(42,4).Value = '8.154'
$ws.Cells.Item(42,4).ClearFormats()
$ws.Cells.Item(42,5).Value = '  -10.87%  '
$ws.Cells.Item(43,4).NumberFormat = '@'
$ws.Cells.Item(43,4).Value = '0.9993'
$ws.Cells.Item(43,4).ClearFormats()
$ws.Cells.Item(43,5).Value = '  +0.01%  '
$ws.Cells.Item(44,4).NumberFormat = '@'
$ws.Cells.Item(44,4).Value = '0.5915'
$ws.Cells.Item(44,4).ClearFormats()
$ws.Cells.Item(44,5).Value = '  -8.16%  '
$ws.Cells.Item(45,4).NumberFormat = '@'
$ws.Cells.Item(45,4).Value = '3.844'
$ws.Cells.Item(45,4).ClearFormats()
$ws.Cells.Item(45,5).Value = '  +0.24%  '
$ws.Cells.Item(46,4).NumberFormat = '@'
$ws.Cells.Item(46,4).Value = '12.82'
$ws.Cells.Item(46,4).ClearFormats()
$ws.Cells.Item(46,5).Value = '  -7.58%  '
$ws.Cells.Item(47,4).NumberFormat = '@'
$ws.Cells.Item(47,4).Value = '0.5623'
$ws.Cells.Item(47,4).ClearFormats()
$ws.Cells.Item(47,5).Value = '  -8.57%  '
$ws.Cells.Item(48,4).NumberFormat = '@'
$ws.Cells.Item(48,4).Value = '124.60'
$ws.Cells.Item(48,4).ClearFormats()
$ws.Cells.Item(48,5).Value = '  -4.00%  '
$ws.Cells.Item(49,4).NumberFormat = '@'
$ws.Cells.Item(49,4).Value = '1.951'
$ws.Cells.Item(49,4).ClearFormats()
$ws.Cells.Item(49,5).Value = '  -7.65%  '
$ws.Cells.Item(50,4).NumberFormat = '@'
$ws.Cells.Item(50,4).Value = '0.06944'
$ws.Cells.Item(50,4).ClearFormats()
$ws.Cells.Item(50,5).Value = '  -5.10%  '
$ws.Cells.Item(51,4).NumberFormat = '@'
$ws.Cells.Item(51,4).Value = '1.192'
$ws.Cells.Item(51,4).ClearFormats()
$ws.Cells.Item(51,5).Value = '  -3.07%  '
